$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 22 (2025-09) stats
$ws.Range("B22").Value = 6296
$ws.Range("D22").Value = 5849514
$ws.Range("E22").Value = 929.0841804320203
$ws.Range("F22").Value = 8.383542778447239
$ws.Range("H22").Value = 27.20870049321933
